# Update "想去人数" (column F) counts on both the "展览" and "全部类型"
# worksheets, which carry duplicate data in this workbook.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1882
    6  = 749
    8  = 359
    9  = 4485
    11 = 350
    12 = 1266
    13 = 528
    15 = 849
    16 = 28
    17 = 464
    19 = 229
    20 = 20
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
